$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new job was reported in Slack - insert a fresh row 2 for it, which
# pushes the previously-existing row 2 (job ending 0785972311) down to row 3.
$ws.Rows.Item(2).Insert()

# Force the new row's cells to text format so numeric-looking values
# (phone numbers / driver numbers with leading zeros) are kept verbatim.
$ws.Range("A2:F2").NumberFormat = "@"

# Populate the new row with the Slack job details.
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = "07777497166"
$ws.Range("C2").Value = "07777497166"
$ws.Range("D2").Value = "046"
$ws.Range("E2").Value = "07777497166 job no 7849325 Pick up at 2 Primrose close, going to office station yard. owes drv 046 £25 Card failed in the car said he will pay in a few days"
$ws.Range("F2").Value = "U092FMBAUP7"

# The row-insert copies the bold header style down onto row 2; the source
# workbook keeps this data row style-free (matching row 3 / the rest of the
# data rows), so strip the inherited formatting back off.
$ws.Range("A2:F2").ClearFormats()
